# Weekly update: add two new daily price records ("Fruta / hortaliza, semanal")
# to the "Hortaliza, Femacal de La Calera - Ajo" sheet.
#
# The existing data block (rows 410-464) holds one record per row; this
# change inserts one new record right before the current row 410, and a
# second new record further down (ending up at row 450 once the first
# insert has shifted everything by one). Inserting real rows (instead of
# overwriting) shifts all the following records down and keeps their
# values/styles intact, growing the used range from A1:R464 to A1:R466.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-AjoRow($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = 3
    $ws.Cells.Item($Row, 2).Value = 'Femacal de La Calera'
    $ws.Cells.Item($Row, 3).Value = 'Coquimbo'
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 5
    $ws.Cells.Item($Row, 6).Value = 100112003
    $ws.Cells.Item($Row, 7).Value = 'Ajo'
    $ws.Cells.Item($Row, 8).Value = 'Chino'
    $ws.Cells.Item($Row, 9).Value = 'Primera'
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = '$/caja 10 kilos'
    $ws.Cells.Item($Row, 15).Value = 'China'
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 10
    $ws.Cells.Item($Row, 18).Value = 'Hortaliza'
}

# Insert the first new record at row 410 (pushes old rows 410-464 down to 411-465).
$ws.Rows.Item(410).Insert()
Set-AjoRow 410 44748 73 16500 17000 16740 1674

# Insert the second new record at row 450 (final numbering), pushing the
# remaining tail down one more row so it ends at 466.
$ws.Rows.Item(450).Insert()
Set-AjoRow 450 44747 50 17000 17000 17000 1700
